$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeDataBatch16")

# Update the photograph paths for rows 3 and 4 (row 2 keeps its original photo)
$ws.Range("D3").Value = "C:\Users\magre\IdeaProjects\HRMSB16\src\test\resources\testdata\Batch16_1.jpg"
$ws.Range("D4").Value = "C:\Users\magre\IdeaProjects\HRMSB16\src\test\resources\testdata\Batch16_2.jpg"

# Update the generated usernames for all three data rows
$ws.Range("E2").Value = "donc1239"
$ws.Range("E3").Value = "donc4569"
$ws.Range("E4").Value = "donc7899"

# Widen column D to fit the longer file paths
$ws.Columns.Item(4).ColumnWidth = 72
